$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework Data Model")

# Update tooltip / component / options for the "Company exchange status" row (row 5)
$ws.Range("F5").Value = "`"Listed`" if the company is listed on an exchange - otherwise `"Unlisted`""
$ws.Range("H5").Value = "Single-Select Radio Button"
$ws.Range("I5").Value = "Listed | Unlisted"

$ws.Range("G1").Select()
$ws.Range("J5").Select()
